$d = $word.ActiveDocument

# Locate the bookmark paragraph (the one holding the _GoBack bookmark),
# which is always the last paragraph in the body. The new paragraph of
# text must be inserted right before it (after the two empty paragraphs
# that follow the picture).
$lastParaIndex = $d.Paragraphs.Count
$bookmarkPara = $d.Paragraphs.Item($lastParaIndex)

# Create a brand-new empty paragraph right before it.
$newParaRange = $bookmarkPara.Range.InsertParagraphBefore()

# The freshly inserted paragraph is now the one before the bookmark paragraph.
$targetPara = $d.Paragraphs.Item($lastParaIndex)
$targetRange = $targetPara.Range

$run1 = "Testul de performanta ajuta la identificarea problemelor aparute in dezvoltarea unui site. De asemenea cu ajutorul testului de performanta putem optimiza aplicatia pentru a fi mai rapida si mai eficienta pentru utilizatori. In graficul de mai sus s-au masurat caracteristici referitoare la useri care au accesat pagina web si timpul de raspuns de la server. Putem observa ca s-au facut un total de 2303 request-uri"
$run2 = " cu media de 13 request-uri pe secunda. Dupa parerea mea rezultatul este destul de bun dar poate fi imbunatatit prin refactorizarea si optimizarea codului aplicatiei."

# Build a minimal WordprocessingML package fragment holding two distinct
# <w:r> runs so they are NOT coalesced into a single run on save (plain
# Range.Text/InsertAfter merges adjacent same-format runs).
$xmlFrag = '<?xml version="1.0" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p>' +
              '<w:r><w:t>' + $run1 + '</w:t></w:r>' +
              '<w:r><w:t xml:space="preserve">' + $run2 + '</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$targetRange.InsertXML($xmlFrag)

Write-Output "done"
